$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("coisas a fazer")

# Row 8: highlight the "Quant. de crismandos" task (green fill, style already used
# elsewhere in the sheet) and add the date it was picked up.
$ws.Range("B8").Interior.Color = 5287936

[void]$ws.Range("C2").Copy()
[void]$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C8").Value = 45784

# New row 10: next item in the to-do list, styled like the row above it (row 9).
[void]$ws.Range("A9").Copy()
[void]$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = 9

[void]$ws.Range("B9").Copy()
[void]$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = "Aprender a fazer backup do banco de dados"

# Update the active-cell selection stored in the sheet view.
[void]$ws.Activate()
[void]$ws.Range("B6").Select()
